$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.790.25"
$ws.Range("E2").Value = "  -7.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.694.76"
$ws.Range("E3").Value = "  -6.62%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.74"
$ws.Range("E5").Value = "  -5.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.52"
$ws.Range("E6").Value = "  +6.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.686.29"
$ws.Range("E7").Value = "  -6.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.631"
$ws.Range("E8").Value = "  -6.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.712"
$ws.Range("E10").Value = "  -5.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("E11").Value = "  -10.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.48"
$ws.Range("E12").Value = "  -5.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  -9.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.59"
$ws.Range("E14").Value = "  -4.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.289.80"
$ws.Range("E15").Value = "  -6.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.744.92"
$ws.Range("E16").Value = "  -5.35%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.127"
$ws.Range("E17").Value = "  -3.28%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.31"
$ws.Range("E18").Value = "  -5.53%  "

$ws.Range("E19").Value = "  -8.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.94"
$ws.Range("E20").Value = "  -7.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.791.71"
$ws.Range("E21").Value = "  -6.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.05"
$ws.Range("E22").Value = "  -7.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.57"
$ws.Range("E23").Value = "  -6.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.20"
$ws.Range("E24").Value = "  -7.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("E25").Value = "  -8.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.80"
$ws.Range("E26").Value = "  -9.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.70"
$ws.Range("E27").Value = "  -2.93%  "

$ws.Range("E28").Value = "  -5.32%  "

$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.50"
$ws.Range("E30").Value = "  -8.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.85"
$ws.Range("E32").Value = "  -8.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.69"
$ws.Range("E33").Value = "  -6.72%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.118"
$ws.Range("E34").Value = "  -8.92%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.33"
$ws.Range("E35").Value = "  -6.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.54"
$ws.Range("E36").Value = "  -6.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0916"
$ws.Range("E37").Value = "  -10.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "597.32"
$ws.Range("E38").Value = "  -6.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  -6.93%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.32"
$ws.Range("E41").Value = "  +14.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -6.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.07"
$ws.Range("E44").Value = "  -11.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0439"
$ws.Range("E45").Value = "  -8.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.46"
$ws.Range("E46").Value = "  -11.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("E48").Value = "  -9.16%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.70"
$ws.Range("E49").Value = "  -14.00%  "

$ws.Range("E50").Value = "  -7.79%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.741.43"
$ws.Range("E51").Value = "  -3.63%  "
